$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing table (symbol,
# timeframe, ...) one column to the right to make room for a new "goal"
# column.
$ws.Columns("A:A").Insert()

# New header + first data row for the "goal" / "best correl-thresh" column.
$ws.Range("A1").Value = "goal"
$ws.Range("A2").Value = "best correl-thresh"

# Widen the new column to fit its contents (matches the width used for the
# other auto-fit columns already on the sheet). Excel's ColumnWidth property
# is in "characters"; the stored sheet width is characters + 5/6, so backing
# out 5/6 here lands the saved width on exactly 17.
$ws.Columns("A:A").ColumnWidth = 17 - (5/6)

# New data row under the (now shifted) header row.
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 93
$ws.Range("G2").Value = 0.4
$ws.Range("I2").Value = "ac"
$ws.Range("J2").Value = 3000
$ws.Range("L2").Value = 5400
$ws.Range("M2").Value = -700
$ws.Range("N2").Value = 0.99

# Match the author's final selection/scroll position.
$ws.Range("N3").Select()
